$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the strain name in cell B3 from "W6t" to "W6"
$ws.Range("B3").Value = "W6"

# Update the active selection to B4, matching the author's final cursor position
$ws.Range("B4").Select()
